# Edit script: insert 4 new data rows (for "Tuna" price records) right before the
# existing row 412, shifting the former rows 412:442 down to 416:446, and then
# populate the newly inserted rows 412:415 with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 412 (pushes old 412..442 down to 416..446)
$ws.Rows("412:415").Insert()

# Common values shared by the 4 new rows
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$fecha       = 45013
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100107
$producto    = "Otros"
$categoriaId = 100107011
$categoria   = "Tuna"
$variedad    = "Sin especificar"
$unidad      = "$/caja 18 kilos"
$kgPorUnidad = 18

# Row 412: Especial
$r = 412
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 250
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 14000
$ws.Cells.Item($r, 16).Value = 14000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 778
$ws.Cells.Item($r, 20).Value = $kgPorUnidad

# Row 413: Extra (doble especial)
$r = 413
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Extra (doble especial)"
$ws.Cells.Item($r, 13).Value = 270
$ws.Cells.Item($r, 14).Value = 16000
$ws.Cells.Item($r, 15).Value = 16000
$ws.Cells.Item($r, 16).Value = 16000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 889
$ws.Cells.Item($r, 20).Value = $kgPorUnidad

# Row 414: Primera
$r = 414
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 280
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 667
$ws.Cells.Item($r, 20).Value = $kgPorUnidad

# Row 415: Segunda
$r = 415
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 556
$ws.Cells.Item($r, 20).Value = $kgPorUnidad

Write-Host "Done inserting rows 412-415"
